# Update the "取得日時" (acquired timestamp) column (A) for all data rows
# on the "ランサーズ" sheet from 2025-09-17 18:25:35 to 2025-09-17 18:38:12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-17 18:25:35"
$newValue = "2025-09-17 18:38:12"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
